# act tablas web jul25
$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsMeta = $wb.Worksheets.Item("Metadata")

# ------------------------------------------------------------------
# Data sheet: insert 3 new rows (2024, 2023, 2022) right after the
# header row, pushing the existing year rows down.
# ------------------------------------------------------------------
$wsData.Rows.Item(2).Resize(3).Insert()

# Force the year labels to be stored as text (matching the existing
# "Fecha" column, where every year is a text value, not a number).
$wsData.Range("A2:A4").NumberFormat = "@"

$wsData.Range("A2").Value = "2024"
$wsData.Range("B2").Value = 37.1

$wsData.Range("A3").Value = "2023"
$wsData.Range("B3").Value = 37.5

$wsData.Range("A4").Value = "2022"
$wsData.Range("B4").Value = 35.5

# Re-apply the plain (unformatted) style used by the rest of the year
# column so the new cells look exactly like their neighbours.
$wsData.Range("A5").Copy()
$wsData.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Metadata sheet: update "observaciones" and "cita" values, and
# insert a new "actualizacion" / "Julio 2025" row right before "cita".
# ------------------------------------------------------------------
$wsMeta.Range("A1").Value = " "

$wsMeta.Range("B8").Value = "Desde marzo de 2020 hasta junio de 2021 se interrumpió el relevamiento presencial y se aplicó de manera telefónica un cuestionario restringido con el objetivo de continuar publicando los indicadores de ingresos y mercado de trabajo. En ese período la encuesta pasó a ser de paneles rotativos elegidos al azar a partir de los casos respondentes del año anterior. `nEn julio de 2021 el INE retomó la realización de encuestas presenciales, pero introdujo un cambio metodológico, ya que la ECH pasa a ser una encuesta de panel rotativo con periodicidad mensual compuesta por seis paneles o grupos de rotación, cada uno de los cuales es una muestra representativa de la población. Con esta nueva metodología, cada hogar seleccionado participa durante seis meses de la ECH."

$wsMeta.Rows.Item(9).Insert()
$wsMeta.Range("A9").Value = "actualizacion"
$wsMeta.Range("B9").Value = "Julio 2025"

$wsMeta.Range("A10").Value = "cita"
$wsMeta.Range("B10").Value = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE (Hasta 2019) / A partir de 2020 con base en ECH - INE`n"
